$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 539.7372594187245
$ws.Range("D2").Value = 120.6474353692902
$ws.Range("G2").Value = 503
$ws.Range("H2").Value = 592
$ws.Range("C3").Value = 44.17396072102088
$ws.Range("D3").Value = 4.758166755696192
$ws.Range("F3").Value = 40.77
$ws.Range("G3").Value = 44.18
$ws.Range("H3").Value = 47.46
$ws.Range("C4").Value = 1.719347887877447
$ws.Range("D4").Value = 3.034312183912127
$ws.Range("F4").Value = 0.65
$ws.Range("G4").Value = 1.33
$ws.Range("H4").Value = 2.32
$ws.Range("C5").Value = 324.1342573648086
$ws.Range("D5").Value = 9.839614067400744
$ws.Range("F5").Value = 319.48
$ws.Range("G5").Value = 325.9
$ws.Range("H5").Value = 331.43
$ws.Range("C6").Value = 22.32895665520659
$ws.Range("D6").Value = 1.905973564937634
$ws.Range("F6").Value = 21.26
$ws.Range("G6").Value = 22.1
$ws.Range("H6").Value = 23
$ws.Range("C7").Value = -76.54445825805192
$ws.Range("D7").Value = 23.67145994587697
$ws.Range("E7").Value = -128
$ws.Range("F7").Value = -93
$ws.Range("C8").Value = 7.603666585792167
$ws.Range("D8").Value = 6.778738504416467
$ws.Range("C9").Value = 9.323404013415951
$ws.Range("D9").Value = 1.688146978757465
$ws.Range("C10").Value = 867.8305786308711
$ws.Range("D10").Value = 0.4612136517686191
$ws.Range("C11").Value = 0.5570904454402326
$ws.Range("D11").Value = 0.590697459462397
$ws.Range("C12").Value = 23.89064813281495
$ws.Range("D12").Value = 13.40902774681995
$ws.Range("C13").Value = 0.6728366531304298
$ws.Range("D13").Value = 0.7487252749585708
$ws.Range("C14").Value = 1.82979191434355
$ws.Range("D14").Value = 1.668613383530017
$ws.Range("C15").Value = 93.94445825805194
$ws.Range("D15").Value = 23.67145994588451
$ws.Range("H15").Value = 110.4
$ws.Range("I15").Value = 145.4
$ws.Range("C16").Value = -85.68788297974999
$ws.Range("D16").Value = 21.36988236077568
$ws.Range("E16").Value = -136.4668316388797
$ws.Range("F16").Value = -102.7643486243649
$ws.Range("G16").Value = -84.59612087980607
$ws.Range("H16").Value = -66.0778545523916
$ws.Range("C17").Value = -78.08421639395782
$ws.Range("D17").Value = 25.88602677945289
$ws.Range("E17").Value = -145.8227324995171
$ws.Range("F17").Value = -93.3707776445072
$ws.Range("G17").Value = -74.1773721860196
$ws.Range("H17").Value = -55.26572375596102
